# Apply cryptocurrency list updates (prices, volumes) and the
# Stacks/SuiNetwork row swap described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'59.452.60"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.78%  "

$ws.Cells.Item(3, 4).Value = "'2.600.64"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.51%  "

$ws.Cells.Item(4, 5).Value = "  +0.00%  "

$ws.Cells.Item(5, 4).Value = "'537.21"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.80%  "

$ws.Cells.Item(6, 4).Value = "'141.34"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.35%  "

$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.18%  "

$ws.Cells.Item(8, 4).Value = "'0.566"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.13%  "

$ws.Cells.Item(9, 4).Value = "'6.49"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.98%  "

$ws.Cells.Item(10, 5).Value = "  +1.29%  "

$ws.Cells.Item(11, 5).Value = "  +1.57%  "

$ws.Cells.Item(12, 5).Value = "  -1.03%  "

$ws.Cells.Item(13, 4).Value = "'3.058.34"
$ws.Cells.Item(13, 4).Style = "Normal"

$ws.Cells.Item(14, 4).Value = "'59.381.17"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.75%  "

$ws.Cells.Item(15, 4).Value = "'20.71"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.94%  "

$ws.Cells.Item(16, 4).Value = "'2.607.52"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.13%  "

$ws.Cells.Item(17, 5).Value = "  +0.23%  "

$ws.Cells.Item(18, 4).Value = "'341.33"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.83%  "

$ws.Cells.Item(19, 4).Value = "'4.37"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.63%  "

$ws.Cells.Item(20, 4).Value = "'10.08"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.15%  "

$ws.Cells.Item(21, 4).Value = "'6.36"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -2.22%  "

$ws.Cells.Item(22, 5).Value = "  +0.01%  "

$ws.Cells.Item(23, 4).Value = "'67.53"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +2.11%  "

$ws.Cells.Item(24, 4).Value = "'0.408"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.17%  "

$ws.Cells.Item(25, 5).Value = "  -1.68%  "

$ws.Cells.Item(26, 4).Value = "'0.999"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.09%  "

$ws.Cells.Item(27, 4).Value = "'7.23"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.87%  "

$ws.Cells.Item(28, 4).Value = "'0.0₃0743"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +2.43%  "

$ws.Cells.Item(29, 5).Value = "  +0.02%  "

$ws.Cells.Item(30, 4).Value = "'1.66"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +5.07%  "

$ws.Cells.Item(31, 5).Value = "  -1.44%  "

$ws.Cells.Item(32, 4).Value = "'18.80"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.55%  "

$ws.Cells.Item(33, 4).Value = "'150.01"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.60%  "

$ws.Cells.Item(34, 4).Value = "'3.97"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.40%  "

$ws.Cells.Item(35, 4).Value = "'1.12"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.33%  "

$ws.Cells.Item(36, 2).Value = "SuiNetwork"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(36, 4).Value = "'0.837"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +2.39%  "

$ws.Cells.Item(37, 2).Value = "Stacks"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(37, 4).Value = "'1.46"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.51%  "

$ws.Cells.Item(38, 4).Value = "'0.824"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.43%  "

$ws.Cells.Item(39, 4).Value = "'3.53"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.62%  "

$ws.Cells.Item(40, 5).Value = "  +0.29%  "

$ws.Cells.Item(41, 4).Value = "'271.68"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.30%  "

$ws.Cells.Item(42, 5).Value = "  +1.80%  "

$ws.Cells.Item(43, 4).Value = "'10.73"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.39%  "

$ws.Cells.Item(44, 5).Value = "  -0.14%  "

$ws.Cells.Item(45, 5).Value = "  +1.49%  "

$ws.Cells.Item(46, 4).Value = "'18.60"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +3.51%  "

$ws.Cells.Item(47, 4).Value = "'1.941.79"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.09%  "

$ws.Cells.Item(48, 5).Value = "  +1.09%  "

$ws.Cells.Item(49, 4).Value = "'4.49"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.33%  "

$ws.Cells.Item(50, 4).Value = "'111.12"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.15%  "

$ws.Cells.Item(51, 4).Value = "'4.77"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.81%  "
